# Apply the update described by the diff:
#  - Remove the 2000-2009 rows (old rows 2-11), shifting the remaining
#    rows (2010-2020) up to rows 2-12.
#  - Update the 2020 row values to the new higher-precision figures.
#  - Append two new rows for 2021 and 2022.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows for years 2000年..2009年 (rows 2-11), shifting rows up.
$ws.Rows("2:11").Delete()

# After the shift, the 2010年..2019年 rows now occupy rows 2-11 unchanged,
# and the former 2020年 row (22) is now row 12. Refresh it with the more
# precise figures from the new source data.
$ws.Range("A12").Value = "2020年"
$ws.Range("B12").Value = 2.23863835626837
$ws.Range("C12").Value = 0.232702784245946
$ws.Range("D12").Value = 1.03609211804918
$ws.Range("E12").Value = 0.969843453973227

# Add the new 2021年 row.
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 8.4
$ws.Range("C13").Value = 0.5
$ws.Range("D13").Value = 4.6
$ws.Range("E13").Value = 3.3

# Add the new 2022年 row.
$ws.Range("A14").Value = "2022年"
$ws.Range("B14").Value = 3
$ws.Range("C14").Value = 0.3
$ws.Range("D14").Value = 1.3
$ws.Range("E14").Value = 1.4

# Copy the style used by the other "A" column year cells (e.g. A2) onto
# the two newly-added year cells so formatting stays consistent.
$ws.Range("A2").Copy()
$ws.Range("A13:A14").PasteSpecial(-4122)

# Make sure the sheet's used-range dimension reflects the new extent.
$ws.Range("A1:E14").Select()
